$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08184"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.76%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.011"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.58%  "

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.819.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.653.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.94%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.71%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3891"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3814"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.34"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.98%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.349"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.19%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.05%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08476"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.99"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.052"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.68%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.094"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.03%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001313"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.93%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.654.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.30%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.11"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07005"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.24%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.61"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.99%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.009"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.21%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.75"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.885.65"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.15%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.433"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.947"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.22%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.93%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.435"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "138.04"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.54%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.841"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.492"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.25%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.836.94"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.30%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02909"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.35%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.662"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.60%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.78"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2676"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09164"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7584"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.88%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.51"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.424"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.46%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.47"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.49%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6938"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.58%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.453"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.116"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.37%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9999"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08285"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.73%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.83"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.226"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.48%  "
